$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("Z1").Value = "Hello"
$v = $ws1.Range("Z1").Value()
Write-Host "Z1 value: [$v]"
